$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("GlobalConstantIntTable")

$ws.Range("A17").Value = "MaxAnalysisLevel"
$ws.Range("B17").Value = 110
